$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5:289 shift down to 6:290.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with a new data point (same constant
# metadata columns as the rest of the table, new date + volume values).
$ws.Cells.Item(5, 1).Value2 = 3
$ws.Cells.Item(5, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(5, 3).Value2 = "Coquimbo"
$ws.Cells.Item(5, 4).Value2 = 44643
$ws.Cells.Item(5, 5).Value2 = 5
$ws.Cells.Item(5, 6).Value2 = 100112039
$ws.Cells.Item(5, 7).Value2 = "Ciboulette"
$ws.Cells.Item(5, 8).Value2 = "Sin especificar"
$ws.Cells.Item(5, 9).Value2 = "Primera"
$ws.Cells.Item(5, 10).Value2 = 120
$ws.Cells.Item(5, 11).Value2 = 1500
$ws.Cells.Item(5, 12).Value2 = 1500
$ws.Cells.Item(5, 13).Value2 = 1500
$ws.Cells.Item(5, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(5, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(5, 16).Value2 = 500
$ws.Cells.Item(5, 17).Value2 = 3
$ws.Cells.Item(5, 18).Value2 = "Hortaliza"
